$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '3323 Villa Maria Catholic Homes St Bernadette''sAged Care Sunshine North'
$ws.Cells.Item(2, 2).Value = 14
$ws.Cells.Item(3, 1).Value = '3398 BlueCross Elly Kay Mordialloc'
$ws.Cells.Item(3, 2).Value = 34
$ws.Cells.Item(4, 1).Value = '3601 Baptcare Westhaven community'
$ws.Cells.Item(4, 2).Value = 20
$ws.Cells.Item(5, 1).Value = '3653 Fronditha Thalpori St Albans Aged Care'
$ws.Cells.Item(5, 2).Value = 22
$ws.Cells.Item(6, 1).Value = '3939 Bupa Aged Care Eastwood'
$ws.Cells.Item(6, 2).Value = 14
$ws.Cells.Item(7, 1).Value = '3975 Aurrum Aged Care Brunswick West'
$ws.Cells.Item(7, 2).Value = 10
$ws.Cells.Item(8, 1).Value = '3988 Kerala Manor Aged Care Diamond Creek'
$ws.Cells.Item(8, 2).Value = 10
$ws.Cells.Item(9, 1).Value = '4257 BlueCross The Gables Camberwell'
$ws.Cells.Item(9, 2).Value = 27
$ws.Cells.Item(10, 1).Value = '4295 Hope Aged Care Sunshine West'
$ws.Cells.Item(10, 2).Value = 31
$ws.Cells.Item(11, 1).Value = '44087 Fitzroy Primary School Fitzroy'
$ws.Cells.Item(11, 2).Value = 22
$ws.Cells.Item(12, 1).Value = '44098 Stawell Primary School'
$ws.Cells.Item(12, 2).Value = 22
$ws.Cells.Item(13, 1).Value = '44234 Lucknow Primary School Bairnsdale'
$ws.Cells.Item(13, 2).Value = 15
$ws.Cells.Item(14, 1).Value = '44366 Lysterfield Primary School Lysterfield'
$ws.Cells.Item(14, 2).Value = 10
$ws.Cells.Item(15, 1).Value = '44444 Nar Nar Goon Primary School Nar NarGoon'
$ws.Cells.Item(15, 2).Value = 18
$ws.Cells.Item(16, 1).Value = '44630 Black Rock Primary School Black Rock'
$ws.Cells.Item(16, 2).Value = 21
$ws.Cells.Item(17, 1).Value = '44811 Dandenong North Primary SchoolDandenong'
$ws.Cells.Item(17, 2).Value = 20
$ws.Cells.Item(18, 1).Value = '44812 Bairnsdale West Primary School'
$ws.Cells.Item(18, 2).Value = 11
$ws.Cells.Item(19, 1).Value = '44865 Parktone Primary School Parkdale'
$ws.Cells.Item(19, 2).Value = 22
$ws.Cells.Item(20, 1).Value = '44950 Templestowe Valley Primary SchoolTemplestowe Lower'
$ws.Cells.Item(20, 2).Value = 27
$ws.Cells.Item(21, 1).Value = '44982 Diamond Creek East Primary SchoolDiamond Creek'
$ws.Cells.Item(21, 2).Value = 10
$ws.Cells.Item(22, 1).Value = '45248 Brookside P-9 College Caroline Springs'
$ws.Cells.Item(22, 2).Value = 30
$ws.Cells.Item(23, 1).Value = '45249 Creekside K-9 College Caroline Springs'
$ws.Cells.Item(23, 2).Value = 12
$ws.Cells.Item(24, 1).Value = '45267 Epping Views Primary School Epping'
$ws.Cells.Item(24, 2).Value = 20
$ws.Cells.Item(25, 1).Value = '45315 Red Hill Consolidated School Red Hill'
$ws.Cells.Item(25, 2).Value = 11
$ws.Cells.Item(26, 1).Value = '45518 Ashwood High School Ashwood'
$ws.Cells.Item(26, 2).Value = 21
$ws.Cells.Item(27, 1).Value = '45569 Nhill College Nhill'
$ws.Cells.Item(27, 2).Value = 33
$ws.Cells.Item(28, 1).Value = '45585 Mount Ridley College Craigieburn'
$ws.Cells.Item(28, 2).Value = 10
$ws.Cells.Item(29, 1).Value = '45648 St Brendans Primary School Shepparton'
$ws.Cells.Item(29, 2).Value = 27
$ws.Cells.Item(30, 1).Value = '4574 Village Glen Aged Care ResidencesMornington'
$ws.Cells.Item(30, 2).Value = 11
$ws.Cells.Item(31, 1).Value = '45784 Holy Rosary Primary School White Hills'
$ws.Cells.Item(31, 2).Value = 26
$ws.Cells.Item(32, 1).Value = '45846 St Mary''s School Mooroopna'
$ws.Cells.Item(32, 2).Value = 15
$ws.Cells.Item(33, 1).Value = '45848 St Kevin''s College Toorak OutbreakGlendalough Campus Junior School'
$ws.Cells.Item(33, 2).Value = 16
$ws.Cells.Item(34, 1).Value = '45950 St. Luke Primary School Lalor'
$ws.Cells.Item(34, 2).Value = 15
$ws.Cells.Item(35, 1).Value = '46028 St Anne''s Catholic Primary SchoolSunbury'
$ws.Cells.Item(35, 2).Value = 11
$ws.Cells.Item(36, 1).Value = '46037 Nazareth Catholic Primary SchoolGrovedal'
$ws.Cells.Item(36, 2).Value = 27
$ws.Cells.Item(37, 1).Value = '46050 Our Lady''s Catholic Primary SchoolCraigieburn'
$ws.Cells.Item(37, 2).Value = 11
$ws.Cells.Item(38, 1).Value = '46052 St. Francis of Assisi Primary School MillPark'
$ws.Cells.Item(38, 2).Value = 26
$ws.Cells.Item(39, 1).Value = '46093 St Brendan''s Primary School Somerville'
$ws.Cells.Item(39, 2).Value = 14
$ws.Cells.Item(40, 1).Value = '46095 Bethany Catholic Primary SchoolWerribee'
$ws.Cells.Item(40, 2).Value = 11
$ws.Cells.Item(41, 1).Value = '46105 Christ the Priest Primary School CarolineSprings'
$ws.Cells.Item(41, 2).Value = 39
$ws.Cells.Item(42, 1).Value = '46125 Our Lady of the Southern Cross PrimarySchool Manor Lakes'
$ws.Cells.Item(42, 2).Value = 36
$ws.Cells.Item(43, 1).Value = '46239 Gilson College Taylors Hill'
$ws.Cells.Item(43, 2).Value = 12
$ws.Cells.Item(44, 1).Value = '46390 Al Siraat College Epping'
$ws.Cells.Item(44, 2).Value = 30
$ws.Cells.Item(45, 1).Value = '50584 St Mary of the Cross MacKillop PrimarySchool Epping'
$ws.Cells.Item(45, 2).Value = 13
$ws.Cells.Item(46, 1).Value = '51529 Sirius College Primary School Dallas'
$ws.Cells.Item(46, 2).Value = 13
$ws.Cells.Item(47, 1).Value = 'Alfred Health The Alfred Hospital Melbourne'
$ws.Cells.Item(47, 2).Value = 13
$ws.Cells.Item(48, 1).Value = 'Camp Coolamatong Farm Camp BanksiaPeninsula'
$ws.Cells.Item(48, 2).Value = 11
$ws.Cells.Item(49, 1).Value = 'Churchill North Primary School Churchill'
$ws.Cells.Item(49, 2).Value = 10
$ws.Cells.Item(50, 1).Value = 'Covenant College Bell Post Hill'
$ws.Cells.Item(50, 2).Value = 23
$ws.Cells.Item(51, 1).Value = 'Epping Views Primary School Camp CapeSchanck'
$ws.Cells.Item(51, 2).Value = 13
$ws.Cells.Item(52, 1).Value = 'Hamilton Country Music Festival Hamilton GolfClub Hamilton'
$ws.Cells.Item(52, 2).Value = 12
$ws.Cells.Item(53, 1).Value = 'Islamic College of Melbourne Tarneit Oct Nov'
$ws.Cells.Item(53, 2).Value = 16
$ws.Cells.Item(54, 1).Value = 'Little Munchkins Childcare Centre Hillside'
$ws.Cells.Item(54, 2).Value = 10
$ws.Cells.Item(55, 1).Value = 'Oakleigh Grammar Melbourne Private SchoolOakleigh'
$ws.Cells.Item(55, 2).Value = 24
$ws.Cells.Item(56, 1).Value = 'Social Gathering 20 November Sunbury'
$ws.Cells.Item(56, 2).Value = 20
$ws.Cells.Item(57, 1).Value = 'Springside Primary School Caroline Springs Nov'
$ws.Cells.Item(57, 2).Value = 23
$ws.Cells.Item(58, 1).Value = 'St Josephs Catholic Primary School Warragul'
$ws.Cells.Item(58, 2).Value = 15
$ws.Cells.Item(59, 1).Value = 'Wagstaff Meat Processing Plant CranbourneEast'
$ws.Cells.Item(59, 2).Value = 36
$ws.Cells.Item(60, 1).Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Cells.Item(60, 2).Value = 14
$ws.Cells.Item(61, 1).Value = 'Western Health Sunshine Hospital EmergencyDepartment St Albans'
$ws.Cells.Item(61, 2).Value = 10
